$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix part orientations for production
$ws.Range("E7").Value = -90
$ws.Range("E8").Value = 90
$ws.Range("E9").Value = 90

# Update the active cell/selection to E9
$ws.Range("E9").Select()
